# Apply the "readd common survey translations" edit.
#
# Summary of the change:
#  - The row on "common_translations" that held the orphaned
#    select_group / Please Select Group strings is removed.
#  - Three new rows are (re)added at the bottom of
#    "framework_translations" (barcode button label + numeric /
#    integer validation messages), each with English, Greek and
#    Spanish text, word-wrapped in the token/English columns.
#  - The active sheet/tab moves from "common_translations" to
#    "framework_translations".

$wb = $excel.ActiveWorkbook

$wsCommon    = $wb.Worksheets.Item("common_translations")
$wsFramework = $wb.Worksheets.Item("framework_translations")

# --- common_translations: drop the leftover select_group row -----------
$wsCommon.Rows.Item(61).Delete()

# --- framework_translations: add the new translation rows --------------
$wsFramework.Range("A55").Value = "barcode_button_label"
$wsFramework.Range("B55").Value = "Scan Barcode"
$wsFramework.Range("C55").Value = "Scan Barcode"
$wsFramework.Range("F55").Value = "Escanear Código de Barras"

$wsFramework.Range("A56").Value = "invalid_numeric_message"
$wsFramework.Range("B56").Value = "Numeric value expected"
$wsFramework.Range("C56").Value = "Αναμενόμενη αριθμητική τιμή"
$wsFramework.Range("F56").Value = "Valor numérico esperado"

$wsFramework.Range("A57").Value = "invalid_integer_message"
$wsFramework.Range("B57").Value = "Integer value expected"
$wsFramework.Range("C57").Value = "Αναμενόμενη τιμή ακέραιας"
$wsFramework.Range("F57").Value = "Valor entero esperado"

# token / english columns get word-wrap applied
$wsFramework.Range("A55:B57").WrapText = $true

# --- view state: common_translations loses focus, framework gains it ---
$wsCommon.Activate()
$wsCommon.Range("H41:I41").Select()

$wsFramework.Activate()
$wsFramework.Range("B49").Select()
